# Daily attendance processing - 2026-01-11 14:32:02
# Swap the order of "Recorded By" entries from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" wherever they appear in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colIndex = 7  # Column G = "Recorded By"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colIndex)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
